# Daily automatic weather-data refresh for meteo.cat summary sheet
# Commit: "Update automàtic: dades i banners [2026-02-25 19:50]"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-25 19:48:33'
$ws.Range("H2").Value = "'46%"
$ws.Range("E3").Value = '2026-02-25 19:48:36'
$ws.Range("E4").Value = '2026-02-25 19:48:38'
$ws.Range("O4").Value = '8.5 °C'
$ws.Range("E5").Value = '2026-02-25 19:48:41'
$ws.Range("H5").Value = "'27%"
$ws.Range("N5").Value = '2.6 °C 19:29 TU'
$ws.Range("O5").Value = '5.9 °C'
$ws.Range("E6").Value = '2026-02-25 19:48:44'
$ws.Range("J6").Value = '1021.6 hPa'
$ws.Range("E7").Value = '2026-02-25 19:48:46'
$ws.Range("J7").Value = '1021.2 hPa'
$ws.Range("E8").Value = '2026-02-25 19:48:48'
$ws.Range("J8").Value = '1020.8 hPa'
$ws.Range("O8").Value = '12.0 °C'
$ws.Range("E9").Value = '2026-02-25 19:48:50'
$ws.Range("O9").Value = '10.1 °C'
$ws.Range("E10").Value = '2026-02-25 19:48:53'
$ws.Range("O10").Value = '9.6 °C'
$ws.Range("E11").Value = '2026-02-25 19:48:56'
$ws.Range("O11").Value = '9.1 °C'
$ws.Range("E12").Value = '2026-02-25 19:48:58'
$ws.Range("O12").Value = '10.0 °C'
$ws.Range("E13").Value = '2026-02-25 19:49:00'
$ws.Range("H13").Value = "'64%"
$ws.Range("J13").Value = '1022.6 hPa'
$ws.Range("O13").Value = '6.8 °C'
$ws.Range("E14").Value = '2026-02-25 19:49:02'
$ws.Range("E15").Value = '2026-02-25 19:49:04'
$ws.Range("E16").Value = '2026-02-25 19:49:07'
$ws.Range("H16").Value = "'30%"
$ws.Range("K16").Value = '13.9 MJ/m2'
$ws.Range("N16").Value = '1.3 °C 19:23 TU'
$ws.Range("E17").Value = '2026-02-25 19:49:10'
$ws.Range("N17").Value = '5.9 °C 19:15 TU'
$ws.Range("O17").Value = '9.4 °C'
$ws.Range("E18").Value = '2026-02-25 19:49:12'
$ws.Range("J18").Value = '1021.8 hPa'
$ws.Range("E19").Value = '2026-02-25 19:49:14'
$ws.Range("O19").Value = '12.5 °C'
$ws.Range("E20").Value = '2026-02-25 19:49:17'
$ws.Range("E21").Value = '2026-02-25 19:49:20'
$ws.Range("O21").Value = '10.1 °C'
$ws.Range("E22").Value = '2026-02-25 19:49:22'
$ws.Range("N22").Value = '0.3 °C 19:20 TU'
$ws.Range("O22").Value = '2.6 °C'
$ws.Range("E23").Value = '2026-02-25 19:49:25'
$ws.Range("E24").Value = '2026-02-25 19:49:27'
$ws.Range("J24").Value = '1019.9 hPa'
$ws.Range("E25").Value = '2026-02-25 19:49:29'
$ws.Range("H25").Value = "'33%"
$ws.Range("E26").Value = '2026-02-25 19:49:31'
$ws.Range("J26").Value = '1019.2 hPa'
$ws.Range("N26").Value = '5.9 °C 19:09 TU'
$ws.Range("O26").Value = '10.4 °C'
$ws.Range("E27").Value = '2026-02-25 19:49:33'
$ws.Range("O27").Value = '5.4 °C'
$ws.Range("E28").Value = '2026-02-25 19:49:35'
$ws.Range("E29").Value = '2026-02-25 19:49:38'
$ws.Range("E30").Value = '2026-02-25 19:49:41'
$ws.Range("J30").Value = '1021.7 hPa'
$ws.Range("E31").Value = '2026-02-25 19:49:44'
$ws.Range("J31").Value = '1021.3 hPa'
$ws.Range("E32").Value = '2026-02-25 19:49:46'
$ws.Range("O32").Value = '9.7 °C'
$ws.Range("E33").Value = '2026-02-25 19:49:49'
$ws.Range("E34").Value = '2026-02-25 19:49:51'
$ws.Range("N34").Value = '-0.1 °C 19:29 TU'
$ws.Range("O34").Value = '3.7 °C'
$ws.Range("E35").Value = '2026-02-25 19:49:54'
$ws.Range("J35").Value = '1019.3 hPa'
$ws.Range("E36").Value = '2026-02-25 19:49:57'
$ws.Range("J36").Value = '1021.8 hPa'
$ws.Range("E37").Value = '2026-02-25 19:49:59'
$ws.Range("J37").Value = '1023.2 hPa'
$ws.Range("E38").Value = '2026-02-25 19:50:02'
$ws.Range("E39").Value = '2026-02-25 19:50:05'
$ws.Range("H39").Value = "'44%"
$ws.Range("E40").Value = '2026-02-25 19:50:07'
$ws.Range("E41").Value = '2026-02-25 19:50:09'
$ws.Range("J41").Value = '1020.8 hPa'
$ws.Range("E42").Value = '2026-02-25 19:50:11'
$ws.Range("E43").Value = '2026-02-25 19:50:14'
$ws.Range("E44").Value = '2026-02-25 19:50:17'
$ws.Range("H44").Value = "'44%"
$ws.Range("E45").Value = '2026-02-25 19:50:19'
$ws.Range("J45").Value = '1019.6 hPa'
$ws.Range("O45").Value = '11.1 °C'
$ws.Range("E46").Value = '2026-02-25 19:50:22'
$ws.Range("J46").Value = '1020.6 hPa'
$ws.Range("O46").Value = '9.7 °C'
